$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 - updated values
$ws.Range("A28").Value = 243335
$ws.Range("C28").Value = 40
$ws.Range("D28").Value = 541.0491803278688
$ws.Range("F28").Value = "2025-04-28 12:19:36"
$ws.Range("G28").Value = "2025-04-28 12:19:36"
$ws.Range("H28").Value = "2025-04-29 13:20:39"
$ws.Range("I28").Value = 33004
$ws.Range("L28").Value = 4
$ws.Range("M28").Value = 152
$ws.Range("N28").Value = 39705

# Row 29 - only column N changes
$ws.Range("N29").Value = 39710

# Row 30 - only column N changes
$ws.Range("N30").Value = 39710

# Row 31 - updated values (previously BIMEC 4, now BIMEC 5)
$ws.Range("A31").Value = 251072
$ws.Range("B31").Value = "BIMEC 5"
$ws.Range("C31").Value = 17
$ws.Range("D31").Value = 82.28169014084507
$ws.Range("E31").Value = "2025-04-28 13:49:27"
$ws.Range("F31").Value = "2025-04-28 14:06:27"
$ws.Range("G31").Value = "2025-04-28 14:06:27"
$ws.Range("H31").Value = "2025-04-29 07:28:44"
$ws.Range("I31").Value = 5842
$ws.Range("L31").Value = 3
$ws.Range("N31").Value = 39705

# Row 32 - updated values (previously BIMEC 5, now BIMEC 2)
$ws.Range("A32").Value = 251126
$ws.Range("B32").Value = "BIMEC 2"
$ws.Range("D32").Value = 156.40625
$ws.Range("E32").Value = "2025-04-28 13:50:02"
$ws.Range("F32").Value = "2025-04-28 14:09:02"
$ws.Range("G32").Value = "2025-04-28 14:09:02"
$ws.Range("H32").Value = "2025-04-29 08:45:27"
$ws.Range("I32").Value = 10010
$ws.Range("L32").Value = 4
$ws.Range("N32").Value = 39705

# Remove rows 33-36 (data no longer present, dimension shrinks to A1:N32)
$ws.Rows("33:36").Delete()
